# Add a new TC-004 test-case block ("errorMessage" -> "Enter Username")
# to the "testdata" sheet, following the existing TC-00x layout pattern:
#   - a blank/spacer row
#   - a header row with the test-case id + field name(s)
#   - a data row with the test-case id + expected value(s)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: blank separator row (mirrors rows 3 & 6 before each test case block).
$ws.Range("A9").Value = " "

# Row 10: TC-004 header row (test case id + field name).
$ws.Range("A10").Value = "TC-004"
$ws.Range("B10").Value = "errorMessage"

# Row 11: TC-004 data row (test case id + expected error message).
$ws.Range("A11").Value = "TC-004"
$ws.Range("B11").Value = "Enter Username"

# The expected-value cells in this sheet (D2, D5, B8) use a distinct
# "Comic Sans MS" blue font; copy that exact formatting onto B11 instead of
# re-deriving it property by property (avoids generating duplicate styles).
$ws.Range("D2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the saved cursor/selection to B11, matching the edited workbook.
$ws.Range("B11").Select()
